$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns store plain text values (the source
# sheet uses inline/shared strings, not numbers). Force text storage for the
# whole data range before writing so numeric-looking prices (e.g. "142.53")
# are not silently reinterpreted as numbers, then restore the default style
# so no stray per-cell formatting is left behind.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '63.911.09'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').Value = '3.063.07'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '559.79'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').Value = '142.53'
$ws.Range('E6').Value = '  -2.18%  '
$ws.Range('D8').Value = '3.062.60'
$ws.Range('E8').Value = '  -0.54%  '
$ws.Range('E9').Value = '  +3.49%  '
$ws.Range('E10').Value = '  +0.60%  '
$ws.Range('D11').Value = '6.11'
$ws.Range('E11').Value = '  -5.31%  '
$ws.Range('E12').Value = '  +2.40%  '
$ws.Range('D13').Value = '0.0000232'
$ws.Range('E13').Value = '  +1.39%  '
$ws.Range('D14').Value = '35.36'
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').Value = '3.562.61'
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('D16').Value = '63.926.92'
$ws.Range('E16').Value = '  -0.97%  '
$ws.Range('D17').Value = '3.061.14'
$ws.Range('E17').Value = '  -0.93%  '
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('D19').Value = '6.81'
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').Value = '486.97'
$ws.Range('E20').Value = '  +2.35%  '
$ws.Range('D21').Value = '14.36'
$ws.Range('E21').Value = '  +2.99%  '
$ws.Range('D22').Value = '0.692'
$ws.Range('E22').Value = '  +1.54%  '
$ws.Range('D23').Value = '14.74'
$ws.Range('E23').Value = '  +8.74%  '
$ws.Range('D24').Value = '7.52'
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').Value = '82.51'
$ws.Range('E25').Value = '  +1.71%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  +0.81%  '
$ws.Range('D28').Value = '8.19'
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('D29').Value = '2.06'
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('D31').Value = '26.50'
$ws.Range('E31').Value = '  +1.58%  '
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('D33').Value = '2.57'
$ws.Range('E33').Value = '  +3.37%  '
$ws.Range('D34').Value = '5.75'
$ws.Range('E34').Value = '  +2.80%  '
$ws.Range('D35').Value = '6.28'
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('D36').Value = '54.78'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').Value = '0.0412'
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('D38').Value = '442.23'
$ws.Range('E38').Value = '  -5.79%  '
$ws.Range('D39').Value = '0.0818'
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('D40').Value = '3.044.74'
$ws.Range('E40').Value = '  +2.69%  '
$ws.Range('D41').Value = '8.37'
$ws.Range('E41').Value = '  +1.21%  '
$ws.Range('E42').Value = '  -8.60%  '
$ws.Range('D43').Value = '0.117'
$ws.Range('E43').Value = '  +2.01%  '
$ws.Range('D44').Value = '0.277'
$ws.Range('E44').Value = '  +6.44%  '
$ws.Range('D45').Value = '27.99'
$ws.Range('E45').Value = '  -1.28%  '
$ws.Range('E46').Value = '  +3.94%  '
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('D49').Value = '0.0₃0518'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('D50').Value = '117.54'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('E51').Value = '  +3.14%  '

$dataRange.Style = "Normal"
